# Add CommonName / sci_name columns (J, K) with species common + scientific
# names, mirroring a paste-in from an external source (GitHub-style
# Helvetica styling on the pasted cells), and tidy up the fish-meal/oil
# allocation sheet (krill already present in the source data; this pass
# backfills the CommonName/sci_name lookup columns for every species,
# including krill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (J1:K1) -------------------------------------------------
# Copy formatting from the existing header cell A1 (bordered header style)
# then set the new header text.
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)
$ws.Range("J1").Value = "CommonName"
$ws.Range("K1").Value = "sci_name"

# ---- "Gross energy" helper rows (2:3) -----------------------------------
$ws.Range("J2:K3").Value = "NA"

# ---- Species rows (4:35): CommonName / sci_name pulled in from an
# external (GitHub) source -------------------------------------------------
$ws.Range("J4:J5").Value = "Antarctic krill"
$ws.Range("K4:K5").Value = "Euphausia superba"

$ws.Range("J6:J11").Value = "Atlantic herring"
$ws.Range("K6:K11").Value = "Clupea harengus"

$ws.Range("J12:J13").Value = "Atlantic mackerel"
$ws.Range("K12:K13").Value = "Scomber scombrus"

$ws.Range("J14:J15").Value = "Blue whiting"
$ws.Range("K14:K15").Value = "Micromesistius poutassou"

$ws.Range("J16:J17").Value = "Boarfish"
$ws.Range("K16:K17").Value = "Capros aper"

$ws.Range("J18:J21").Value = "Capelin"
$ws.Range("K18:K21").Value = "Mallotus villosus"

$ws.Range("J22:J23").Value = "Chilean jack"
$ws.Range("K22:K23").Value = "Trachurus murphyi"

$ws.Range("J24:J25").Value = "European sprat"
$ws.Range("K24:K25").Value = "Sprattus sprattus"

$ws.Range("J26:J27").Value = "Gulf menhaden"
$ws.Range("K26:K27").Value = "Brevoortia patronus"

$ws.Range("J28:J29").Value = "Norway pout"
$ws.Range("K28:K29").Value = "Trisopterus esmarkii"

$ws.Range("J30:J31").Value = "Peruvian anchovy"
$ws.Range("K30:K31").Value = "Engraulis ringens"

$ws.Range("J32:J33").Value = "Sandeels"
$ws.Range("K32:K33").Value = "Ammodytes tobianus"

$ws.Range("J34:J35").Value = "South American pilchard"
$ws.Range("K34:K35").Value = "Sardinops sagax"

# Single formatting pass over every pasted species cell (Helvetica 12,
# color #1F2328) -- minimizes spurious intermediate style allocations.
$ws.Range("J4:K35").Font.Name = "Helvetica"
$ws.Range("J4:K35").Font.Size = 12
$ws.Range("J4:K35").Font.Color = 2630431

# ---- Global-average rows (36:37) ----------------------------------------
# Row 37 already carries the table's closing bottom-border; mirror that
# onto K36/J37/K37, and keep J36 border-free (matching the rest of row 36).
$ws.Range("A37").Copy()
$ws.Range("K36").PasteSpecial(-4122)
$ws.Range("J37:K37").PasteSpecial(-4122)

$ws.Range("J36:K37").Value = "NA"

# Single formatting pass for the plain (Calibri, black) "NA" lookup cells.
$ws.Range("J2:K3,J36:K37").Font.Name = "Calibri"
$ws.Range("J2:K3,J36:K37").Font.Size = 11
$ws.Range("J2:K3,J36:K37").Font.Color = 0

# ---- Column widths for the new J/K columns -------------------------------
$ws.Columns("J").ColumnWidth = 25.59
$ws.Columns("K").ColumnWidth = 27.09

# ---- Restore selection (matches the saved file's cursor position) -------
$ws.Range("G9").Select()
